$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed numeric values per diff (data/ifrs/현대로템.xlsx)
$ws.Range("D2").Value = 31911
$ws.Range("E2").Value = 66
$ws.Range("F2").Value = 66
$ws.Range("G2").Value = -115
$ws.Range("H2").Value = -151
$ws.Range("I2").Value = -167
$ws.Range("J2").Value = 16
$ws.Range("K2").Value = 44008
$ws.Range("L2").Value = 26622
$ws.Range("M2").Value = 17387
$ws.Range("N2").Value = 16610
$ws.Range("O2").Value = 777
$ws.Range("P2").Value = 4250
$ws.Range("Q2").Value = -1888
$ws.Range("R2").Value = -578
$ws.Range("S2").Value = 2455
$ws.Range("T2").Value = 406
$ws.Range("U2").Value = -2294
$ws.Range("V2").Value = 13749
$ws.Range("W2").Value = 0.21
$ws.Range("X2").Value = -0.47
$ws.Range("Y2").Value = -0.99
$ws.Range("Z2").Value = -0.35
$ws.Range("AA2").Value = 153.12
$ws.Range("AB2").Value = 290.96
$ws.Range("AC2").Value = -197
$ws.Range("AD2").Value = -101.4
$ws.Range("AE2").Value = 19541
$ws.Range("AF2").Value = 1.02
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 85000000
$ws.Range("D3").Value = 33091
$ws.Range("E3").Value = -1929
$ws.Range("F3").Value = -1929
$ws.Range("G3").Value = -2970
$ws.Range("H3").Value = -3045
$ws.Range("I3").Value = -3062
$ws.Range("J3").Value = 17
$ws.Range("K3").Value = 50439
$ws.Range("L3").Value = 36063
$ws.Range("M3").Value = 14377
$ws.Range("N3").Value = 13587
$ws.Range("O3").Value = 789
$ws.Range("P3").Value = 4250
$ws.Range("Q3").Value = -5645
$ws.Range("R3").Value = -553
$ws.Range("S3").Value = 7642
$ws.Range("T3").Value = 550
$ws.Range("U3").Value = -6194
$ws.Range("V3").Value = 21437
$ws.Range("W3").Value = -5.83
$ws.Range("X3").Value = -9.199999999999999
$ws.Range("Y3").Value = -20.28
$ws.Range("Z3").Value = -6.45
$ws.Range("AA3").Value = 250.84
$ws.Range("AB3").Value = 220.31
$ws.Range("AC3").Value = -3602
$ws.Range("AD3").Value = -4.07
$ws.Range("AE3").Value = 15985
$ws.Range("AF3").Value = 0.92
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 85000000
$ws.Range("D4").Value = 29848
$ws.Range("E4").Value = 1062
$ws.Range("F4").Value = 1062
$ws.Range("G4").Value = 231
$ws.Range("H4").Value = 231
$ws.Range("I4").Value = 218
$ws.Range("J4").Value = 13
$ws.Range("K4").Value = 44732
$ws.Range("L4").Value = 30060
$ws.Range("M4").Value = 14672
$ws.Range("N4").Value = 13920
$ws.Range("O4").Value = 752
$ws.Range("P4").Value = 4250
$ws.Range("Q4").Value = 6111
$ws.Range("R4").Value = 556
$ws.Range("S4").Value = -3435
$ws.Range("T4").Value = 199
$ws.Range("U4").Value = 5912
$ws.Range("V4").Value = 18128
$ws.Range("W4").Value = 3.56
$ws.Range("X4").Value = 0.78
$ws.Range("Y4").Value = 1.59
$ws.Range("Z4").Value = 0.49
$ws.Range("AA4").Value = 204.88
$ws.Range("AB4").Value = 227.66
$ws.Range("AC4").Value = 257
$ws.Range("AD4").Value = 70.66
$ws.Range("AE4").Value = 16376
$ws.Range("AF4").Value = 1.11
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 85000000
$ws.Range("D5").Value = 27257
$ws.Range("E5").Value = 454
$ws.Range("F5").Value = 454
$ws.Range("G5").Value = -479
$ws.Range("H5").Value = -463
$ws.Range("I5").Value = -474
$ws.Range("J5").Value = 11
$ws.Range("K5").Value = 40839
$ws.Range("L5").Value = 26656
$ws.Range("M5").Value = 14183
$ws.Range("N5").Value = 13440
$ws.Range("O5").Value = 743
$ws.Range("P5").Value = 4250
$ws.Range("Q5").Value = 2192
$ws.Range("R5").Value = 312
$ws.Range("S5").Value = -3918
$ws.Range("T5").Value = 280
$ws.Range("U5").Value = 1913
$ws.Range("V5").Value = 14070
$ws.Range("W5").Value = 1.67
$ws.Range("X5").Value = -1.7
$ws.Range("Y5").Value = -3.46
$ws.Range("Z5").Value = -1.08
$ws.Range("AA5").Value = 187.94
$ws.Range("AB5").Value = 216.78
$ws.Range("AC5").Value = -557
$ws.Range("AD5").Value = -33.64
$ws.Range("AE5").Value = 15811
$ws.Range("AF5").Value = 1.19
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 85000000
$ws.Range("D6").Value = 24119
$ws.Range("E6").Value = -1962
$ws.Range("F6").Value = -1962
$ws.Range("G6").Value = -3241
$ws.Range("H6").Value = -3080
$ws.Range("I6").Value = -3008
$ws.Range("K6").Value = 40022
$ws.Range("L6").Value = 28942
$ws.Range("M6").Value = 11080
$ws.Range("N6").Value = 10440
$ws.Range("P6").Value = 4250
$ws.Range("Q6").Value = -142
$ws.Range("R6").Value = -381
$ws.Range("S6").Value = -195
$ws.Range("T6").Value = 300
$ws.Range("U6").Value = -442
$ws.Range("V6").Value = 13921
$ws.Range("W6").Value = -8.130000000000001
$ws.Range("X6").Value = -12.77
$ws.Range("Y6").Value = -25.2
$ws.Range("Z6").Value = -7.62
$ws.Range("AA6").Value = 261.21
$ws.Range("AB6").Value = 146.59
$ws.Range("AC6").Value = -3539
$ws.Range("AD6").Value = -7.87
$ws.Range("AE6").Value = 12282
$ws.Range("AF6").Value = 2.27
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 85000000
$ws.Range("D7").Value = 25388
$ws.Range("E7").Value = -1417
$ws.Range("G7").Value = -1846
$ws.Range("H7").Value = -1831
$ws.Range("I7").Value = -1836
$ws.Range("K7").Value = 40134
$ws.Range("L7").Value = 30511
$ws.Range("M7").Value = 9623
$ws.Range("N7").Value = 8909
$ws.Range("P7").Value = 4250
$ws.Range("Q7").Value = -1210
$ws.Range("R7").Value = 492
$ws.Range("S7").Value = 1626
$ws.Range("T7").Value = 288
$ws.Range("U7").Value = -676
$ws.Range("W7").Value = -5.58
$ws.Range("X7").Value = -7.21
$ws.Range("Y7").Value = -18.98
$ws.Range("Z7").Value = -4.57
$ws.Range("AA7").Value = 317.05
$ws.Range("AC7").Value = -2160
$ws.Range("AD7").Value = -7.52
$ws.Range("AE7").Value = 10481
$ws.Range("AF7").Value = 1.55
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("D8").Value = 27228
$ws.Range("E8").Value = 371
$ws.Range("G8").Value = -65
$ws.Range("H8").Value = -81
$ws.Range("I8").Value = -80
$ws.Range("K8").Value = 41182
$ws.Range("L8").Value = 31999
$ws.Range("M8").Value = 9183
$ws.Range("N8").Value = 8498
$ws.Range("P8").Value = 4250
$ws.Range("Q8").Value = 870
$ws.Range("R8").Value = -345
$ws.Range("S8").Value = 210
$ws.Range("T8").Value = 325
$ws.Range("U8").Value = 583
$ws.Range("W8").Value = 1.36
$ws.Range("X8").Value = -0.3
$ws.Range("Y8").Value = -0.92
$ws.Range("Z8").Value = -0.2
$ws.Range("AA8").Value = 348.48
$ws.Range("AC8").Value = -94
$ws.Range("AD8").Value = -151.14
$ws.Range("AE8").Value = 9998
$ws.Range("AF8").Value = 1.43
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("D9").Value = 29868
$ws.Range("E9").Value = 846
$ws.Range("G9").Value = 452
$ws.Range("H9").Value = 369
$ws.Range("I9").Value = 368
$ws.Range("K9").Value = 41718
$ws.Range("L9").Value = 32267
$ws.Range("M9").Value = 9452
$ws.Range("N9").Value = 8860
$ws.Range("P9").Value = 4250
$ws.Range("Q9").Value = 1052
$ws.Range("R9").Value = -338
$ws.Range("S9").Value = -220
$ws.Range("T9").Value = 348
$ws.Range("U9").Value = 724
$ws.Range("W9").Value = 2.83
$ws.Range("X9").Value = 1.23
$ws.Range("Y9").Value = 4.24
$ws.Range("Z9").Value = 0.89
$ws.Range("AA9").Value = 341.39
$ws.Range("AC9").Value = 433
$ws.Range("AD9").Value = 32.88
$ws.Range("AE9").Value = 10424
$ws.Range("AF9").Value = 1.37
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0
$ws.Range("AI9").Value = 0

# Cells removed entirely in the new version (no longer present)
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI7").ClearContents()
$ws.Range("AI8").ClearContents()
